# Scheduled-runner refresh of market/profit figures across the Sargatanas
# profit-tracking sheets (one sheet per crafting job). Each block below
# updates the currentAveragePrice*/LevePrice*/LeveProfit* columns (H:N) for
# the specific leve rows that moved in this run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2219.7144
$ws.Range("I33").Value = 2907.6
$ws.Range("K33").Value = 2907.6
$ws.Range("M33").Value = -2678.6
$ws.Range("H55").Value = 194.1579
$ws.Range("I55").Value = 199.94444
$ws.Range("K55").Value = 199.94444
$ws.Range("M55").Value = 14.05556000000001
$ws.Range("H132").Value = 2744.9
$ws.Range("I132").Value = 2744.9
$ws.Range("K132").Value = 8234.700000000001
$ws.Range("M132").Value = -5704.700000000001
$ws.Range("H137").Value = 3713.8667
$ws.Range("I137").Value = 2463.9644
$ws.Range("K137").Value = 7391.8932
$ws.Range("M137").Value = -4841.8932
$ws.Range("H138").Value = 2277662.8
$ws.Range("I138").Value = 1863
$ws.Range("J138").Value = 3036262.8
$ws.Range("K138").Value = 5589
$ws.Range("L138").Value = 9108788.399999999
$ws.Range("M138").Value = -449
$ws.Range("N138").Value = -9119068.399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 509
$ws.Range("I25").Value = 509
$ws.Range("K25").Value = 509
$ws.Range("M25").Value = -107
$ws.Range("H32").Value = 3912836.8
$ws.Range("I32").Value = 3912836.8
$ws.Range("K32").Value = 3912836.8
$ws.Range("M32").Value = -3912549.8
$ws.Range("H61").Value = 37043036
$ws.Range("I61").Value = 1481.8125
$ws.Range("K61").Value = 1481.8125
$ws.Range("M61").Value = -1269.8125
$ws.Range("H122").Value = 10887.549
$ws.Range("I122").Value = 17877.785
$ws.Range("J122").Value = 5130.8823
$ws.Range("K122").Value = 53633.355
$ws.Range("L122").Value = 15392.6469
$ws.Range("M122").Value = -51183.355
$ws.Range("N122").Value = -20292.6469
$ws.Range("H136").Value = 37043036
$ws.Range("I136").Value = 1481.8125
$ws.Range("K136").Value = 4445.4375
$ws.Range("M136").Value = -1895.4375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3478207.8
$ws.Range("I107").Value = 4814461.5
$ws.Range("J107").Value = 3948.6
$ws.Range("K107").Value = 4814461.5
$ws.Range("L107").Value = 3948.6
$ws.Range("M107").Value = -4812541.5
$ws.Range("N107").Value = -7788.6
$ws.Range("H134").Value = 5819317.5
$ws.Range("I134").Value = 9617253
$ws.Range("J134").Value = 10711.059
$ws.Range("K134").Value = 28851759
$ws.Range("L134").Value = 32133.177
$ws.Range("M134").Value = -28849224
$ws.Range("N134").Value = -37203.177

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5673.069
$ws.Range("I31").Value = 1997.9286
$ws.Range("K31").Value = 1997.9286
$ws.Range("M31").Value = -1702.9286
$ws.Range("H34").Value = 5673.069
$ws.Range("I34").Value = 1997.9286
$ws.Range("K34").Value = 1997.9286
$ws.Range("M34").Value = -1795.9286
$ws.Range("H58").Value = 7764.6553
$ws.Range("I58").Value = 2870.7856
$ws.Range("K58").Value = 2870.7856
$ws.Range("M58").Value = -2667.7856
$ws.Range("H107").Value = 3349.7778
$ws.Range("I107").Value = 4215
$ws.Range("J107").Value = 2917.1667
$ws.Range("K107").Value = 4215
$ws.Range("L107").Value = 2917.1667
$ws.Range("M107").Value = -2295
$ws.Range("N107").Value = -6757.1667
$ws.Range("H132").Value = 8344.5
$ws.Range("I132").Value = 4303.4
$ws.Range("J132").Value = 10181.363
$ws.Range("K132").Value = 12910.2
$ws.Range("L132").Value = 30544.089
$ws.Range("M132").Value = -10380.2
$ws.Range("N132").Value = -35604.089
$ws.Range("H134").Value = 6564.5
$ws.Range("J134").Value = 8461.958000000001
$ws.Range("L134").Value = 25385.874
$ws.Range("N134").Value = -30455.874
$ws.Range("H136").Value = 7764.6553
$ws.Range("I136").Value = 2870.7856
$ws.Range("K136").Value = 8612.356800000001
$ws.Range("M136").Value = -6062.356800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 4041218.8
$ws.Range("I122").Value = 9428143
$ws.Range("J122").Value = 1026
$ws.Range("K122").Value = 84853287
$ws.Range("L122").Value = 9234
$ws.Range("M122").Value = -84850837
$ws.Range("N122").Value = -14134
$ws.Range("H127").Value = 2857.1428
$ws.Range("J127").Value = 2857.1428
$ws.Range("L127").Value = 8571.428400000001
$ws.Range("N127").Value = -18491.4284
$ws.Range("H132").Value = 6003.1665
$ws.Range("I132").Value = 2005.3684
$ws.Range("K132").Value = 18048.3156
$ws.Range("M132").Value = -15518.3156

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 2000
$ws.Range("I3").Value = 2000
$ws.Range("K3").Value = 2000
$ws.Range("M3").Value = -1884
$ws.Range("H80").Value = 5857.1113
$ws.Range("J80").Value = 7623.75
$ws.Range("L80").Value = 7623.75
$ws.Range("N80").Value = -9619.75
$ws.Range("H83").Value = 5857.1113
$ws.Range("J83").Value = 7623.75
$ws.Range("L83").Value = 38118.75
$ws.Range("N83").Value = -48102.75
$ws.Range("H107").Value = 728031.5600000001
$ws.Range("I107").Value = 1000969.5
$ws.Range("K107").Value = 1000969.5
$ws.Range("M107").Value = -999049.5
$ws.Range("H126").Value = 41675030
$ws.Range("J126").Value = 9127.182000000001
$ws.Range("L126").Value = 27381.546
$ws.Range("N126").Value = -32321.546
$ws.Range("H132").Value = 3652.9756
$ws.Range("I132").Value = 2023.3939
$ws.Range("K132").Value = 6070.1817
$ws.Range("M132").Value = -3540.1817

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3423.25
$ws.Range("I61").Value = 1878.8
$ws.Range("J61").Value = 6933.364
$ws.Range("K61").Value = 1878.8
$ws.Range("L61").Value = 6933.364
$ws.Range("M61").Value = -1676.8
$ws.Range("N61").Value = -7337.364
$ws.Range("H93").Value = 5009.857
$ws.Range("I93").Value = 5015.4
$ws.Range("K93").Value = 5015.4
$ws.Range("M93").Value = -3767.4
$ws.Range("H113").Value = 3423.25
$ws.Range("I113").Value = 1878.8
$ws.Range("J113").Value = 6933.364
$ws.Range("K113").Value = 1878.8
$ws.Range("L113").Value = 6933.364
$ws.Range("M113").Value = 291.2
$ws.Range("N113").Value = -11273.364
$ws.Range("H136").Value = 8476.049000000001
$ws.Range("I136").Value = 3272.75
$ws.Range("K136").Value = 9818.25
$ws.Range("M136").Value = -7268.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 10546823
$ws.Range("I81").Value = 909387
$ws.Range("K81").Value = 1818774
$ws.Range("M81").Value = -1817713
$ws.Range("H84").Value = 10546823
$ws.Range("I84").Value = 909387
$ws.Range("K84").Value = 9093870
$ws.Range("M84").Value = -9088566
$ws.Range("H113").Value = 510.6279
$ws.Range("I113").Value = 476.8095
$ws.Range("J113").Value = 542.9091
$ws.Range("K113").Value = 1430.4285
$ws.Range("L113").Value = 1628.7273
$ws.Range("M113").Value = 739.5715
$ws.Range("N113").Value = -5968.7273
$ws.Range("H126").Value = 3208.3
$ws.Range("I126").Value = 3993.5
$ws.Range("J126").Value = 2684.8333
$ws.Range("K126").Value = 11980.5
$ws.Range("L126").Value = 8054.499899999999
$ws.Range("M126").Value = -9510.5
$ws.Range("N126").Value = -12994.4999
$ws.Range("H136").Value = 42089000
$ws.Range("I136").Value = 100001800
$ws.Range("K136").Value = 300005400
$ws.Range("M136").Value = -300002850
